$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column Y
$ws.Range("Y1").Value = "AuslaenderAnteil"

# AuslaenderAnteil values for rows 2-89 (GemeindeName order as in sheet)
$values = @(
    10.0,
    10.4,
    11.5,
    8.1,
    17.4,
    9.3,
    6.5,
    11.9,
    11.9,
    17.9,
    12.1,
    25.2,
    14.2,
    16.9,
    4.8,
    4.5,
    20.8,
    8.6,
    6.5,
    33.3,
    6.7,
    11.3,
    10.2,
    5.8,
    5.8,
    5.8,
    8.5,
    10.1,
    17.3,
    19.2,
    21.1,
    17.5,
    16.7,
    6.9,
    12.3,
    4.8,
    5.8,
    8.1,
    14.8,
    20.9,
    6.6,
    8.7,
    17.0,
    10.3,
    9.5,
    17.7,
    5.0,
    24.2,
    11.5,
    9.2,
    13.8,
    13.6,
    9.3,
    26.5,
    10.2,
    11.0,
    11.6,
    11.2,
    5.4,
    21.8,
    17.9,
    17.9,
    5.4,
    1.5,
    24.9,
    8.9,
    7.6,
    5.3,
    4.9,
    8.7,
    5.5,
    6.4,
    16.4,
    16.4,
    8.1,
    7.9,
    15.3,
    24.5,
    9.6,
    4.6,
    24.0,
    25.6,
    22.3,
    13.1,
    15.2,
    9.9,
    15.7,
    11.6
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 25).Value = $values[$i]
}
